$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.499.17'
$ws.Range("E2").Value = '  +2.48%  '

# Row 3
$ws.Range("D3").Value = '1.828.62'
$ws.Range("E3").Value = '  +2.03%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.46'
$ws.Range("E5").Value = '  -0.11%  '

# Row 6
$ws.Range("E6").Value = '  +0.11%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5030'
$ws.Range("E7").Value = '  -6.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3910'
$ws.Range("E8").Value = '  +2.24%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07702'
$ws.Range("E9").Value = '  +3.55%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.91'
$ws.Range("E10").Value = '  +1.11%  '

# Row 11
$ws.Range("E11").Value = '  +2.47%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.02'
$ws.Range("E12").Value = '  +3.50%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.249'
$ws.Range("E13").Value = '  +0.83%  '

# Row 14
$ws.Range("E14").Value = '  +0.06%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.555'
$ws.Range("E15").Value = '  +1.72%  '

# Row 16
$ws.Range("D16").Value = '1.825.56'
$ws.Range("E16").Value = '  +2.33%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.43'
$ws.Range("E17").Value = '  +5.82%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001080'
$ws.Range("E18").Value = '  +2.21%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06601'
$ws.Range("E19").Value = '  +1.36%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.73'
$ws.Range("E20").Value = '  +2.63%  '

# Row 21
$ws.Range("E21").Value = '  +0.03%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.136'
$ws.Range("E22").Value = '  +2.96%  '

# Row 23
$ws.Range("D23").Value = '28.529.94'
$ws.Range("E23").Value = '  +2.42%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.13'
$ws.Range("E24").Value = '  +0.13%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.247'
$ws.Range("E25").Value = '  +7.33%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.82'
$ws.Range("E26").Value = '  -0.23%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.63'

# Row 28
$ws.Range("D28").Value = '2.034.78'
$ws.Range("E28").Value = '  +1.82%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.411'
$ws.Range("E29").Value = '  +3.86%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.95'
$ws.Range("E30").Value = '  +2.81%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.136'
$ws.Range("E31").Value = '  +2.87%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1088'
$ws.Range("E32").Value = '  -0.51%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.653'
$ws.Range("E33").Value = '  +2.76%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07068'
$ws.Range("E35").Value = '  +1.78%  '

# Row 36
$ws.Range("E36").Value = '  +1.41%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.958'
$ws.Range("E37").Value = '  +6.55%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02321'
$ws.Range("E38").Value = '  +2.25%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.126'
$ws.Range("E39").Value = '  +1.71%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6236'
$ws.Range("E40").Value = '  +2.37%  '

# Row 41
$ws.Range("E41").Value = '  -1.39%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.193'
$ws.Range("E42").Value = '  +2.37%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.06%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.397'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.43'
$ws.Range("E45").Value = '  +1.16%  '

# Row 46
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.716'
$ws.Range("E46").Value = '  +1.08%  '

# Row 47
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5890'
$ws.Range("E47").Value = '  +3.52%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.28'
$ws.Range("E48").Value = '  -0.62%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.970'
$ws.Range("E49").Value = '  +3.30%  '

# Row 50
$ws.Range("E50").Value = '  +1.15%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06935'
$ws.Range("E51").Value = '  +1.96%  '
